$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1: "Save", copying G1's formatting (bold/centered/bordered header style)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# New data column H2:H3
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
